$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Convert column B (Voltage, mV) into Volts; formulas in D (Power) and new
#     helper column E (old mV-style value, kept as B/1000 for reference) ---

$ws.Range("B3").Value = 0.412
$ws.Range("B4").Value = 0.411
$ws.Range("B5").Value = 0.411
$ws.Range("B6").Value = 0.411
$ws.Range("B7").Value = 0.41
$ws.Range("B8").Value = 0.409
$ws.Range("B9").Value = 0.396
$ws.Range("B10").Value = 0.313
$ws.Range("B11").Value = 0.208
$ws.Range("B12").Value = 0.095
$ws.Range("B13").Value = 0.085
$ws.Range("B14").Value = 0.077

# Power column no longer needs to divide B by 1000 (B is already in volts)
$ws.Range("D3").Formula = "=C3*B3"
$ws.Range("D4:D15").Formula = "=C4*B4"

# New helper column E recreates the original raw (mA-style) reading
$ws.Range("E3").Formula = "=B3/1000"
$ws.Range("E4:E15").Formula = "=B4/1000"

# --- Chart 1 (I-V curve) ---
$chart1 = $ws.ChartObjects().Item(1).Chart
$chart1.HasTitle = $true
$chart1.ChartTitle.Text = "I-V Curve of PV Emulator"

$ivX = $chart1.Axes(1)
$ivX.HasTitle = $true
$ivX.AxisTitle.Text = "Voltage (V)"

$ivY = $chart1.Axes(2)
$ivY.HasTitle = $true
$ivY.AxisTitle.Text = "Current (A)"

# --- Chart 2 (P-V curve) ---
$chart2 = $ws.ChartObjects().Item(2).Chart
$chart2.HasTitle = $true
$chart2.ChartTitle.Text = "P-V Curve of PV Emulator"

$pvX = $chart2.Axes(1)
$pvX.HasTitle = $true
$pvX.AxisTitle.Text = "Voltage (V)"

$pvY = $chart2.Axes(2)
$pvY.HasTitle = $true
$pvY.AxisTitle.Text = "Power (W)"

# --- Selection as left by the author ---
[void]$ws.Range("T25").Select()
